# "analiza wszedzie od 10k" - re-run the timing analysis with measurements
# starting from n=10000 for every series, expressed in microseconds (column B),
# and add a derived "milliseconds" column C = B/1000 that feeds the existing
# TRANSPOSE() summary table (I4:O9) and the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New raw measurements in column B (rows 2:36) ---------------------
$newB = @{
    2  = 14092.7
    3  = 60961.3
    4  = 242347
    5  = 924394
    6  = 3732150
    7  = 14529400
    8  = 58132700
    9  = 10
    10 = 16
    11 = 34
    12 = 76
    13 = 144.667
    14 = 296.667
    15 = 639.333
    16 = 29887.3
    17 = 112000
    18 = 450417
    19 = 1828300
    20 = 7292180
    21 = 29271400
    22 = 116830000
    23 = 13129
    24 = 52681.3
    25 = 204020
    26 = 806882
    27 = 3236500
    28 = 12970500
    29 = 51972600
    30 = 8322.67
    31 = 30627.3
    32 = 124572
    33 = 511414
    34 = 1998730
    35 = 8027180
    36 = 32397600
}

foreach ($row in $newB.Keys) {
    $ws.Range("B$row").Value = $newB[$row]
}

# --- 2. New column C = B/1000 (ms), filled down with a relative formula --
$ws.Range("C2").Formula = "=B2/1000"
$ws.Range("C3:C36").Formula = "=B3/1000"

# --- 3. Repoint the TRANSPOSE() summary table at column C instead of B ---
$ws.Range("I5:O5").FormulaArray = "=TRANSPOSE(C2:C8)"
$ws.Range("I6:O6").FormulaArray = "=TRANSPOSE(C9:C15)"
$ws.Range("I7:O7").FormulaArray = "=TRANSPOSE(C16:C22)"
$ws.Range("I8:O8").FormulaArray = "=TRANSPOSE(C23:C29)"
$ws.Range("I9:O9").FormulaArray = "=TRANSPOSE(C30:C36)"

# --- 4. Selection moved to the first repointed summary row ----------------
$ws.Range("I5:O5").Select()
